$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.490.76'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.51%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.690.23'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +5.77%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '520.25'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +2.52%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.92'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.29%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.994'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.572'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +1.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.727.96'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +7.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.28'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.46%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +5.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.341'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +2.61%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.126'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.161.50'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +5.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '59.449.63'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.48%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.31'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.29%  '
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.00%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.722.51'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +6.93%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '360.29'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +7.54%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.04%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.53'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.30%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.28'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.57%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.16'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.426'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.39%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.162'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.87%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.992'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0825'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.73%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.31'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +5.72%  '
$ws.Range('B30').Value = 'USDe'
$ws.Range('C30').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.996'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.33%  '
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.43'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +9.45%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.27'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.03%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.60'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.73%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '150.80'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('B35').Value = 'SuiNetwork'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.989'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +6.41%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.08'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +5.16%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.16'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +4.68%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '36.99'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +2.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.860'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +4.25%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.75'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +6.33%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.43'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.38%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '283.92'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.626'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.63%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.13'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +7.49%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0992'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.992'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.58%  '
$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0537'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.50%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.79'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.62%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.48%  '
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.018.14'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +6.83%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '10.29'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.01%  '
